$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("accountInfo")

# Add new header columns
$ws.Range("E1").Value = "Friends"
$ws.Range("F1").Value = "Friend Requests"

# Copy the header style (bold) from D1 to the new header cells
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122) # xlPasteFormats

# Add friend request note for Tom (row 14)
$ws.Range("F14").Value = "Bob has sent you a friend request.>>Jose has sent you a friend request."

# Replace row 16 (testnarek) with new Scooby account, clearing the old Library column
$ws.Range("A16").Value = "Scooby"
$ws.Range("B16").Value = '$coobyDoobyD00'
$ws.Range("C16").Value = "scoobysnacks@yahoo.com"
$ws.Range("D16").ClearContents()

# Move selection to E2 to match final cursor position
$ws.Range("E2").Select()
